# Automatic update of files.
# The "Förändrad" (Changed) date column (C) for every data row (2-89)
# advances by one day: 45174 (2023-09-05) -> 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C89").Value = 45175
